# Generate Report for Handback
# Adds a new handback record (afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md) as row 4
# on the Overview sheet, and on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$srcFile   = "afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md"
$srcPath   = "e2e\afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md"
$ext       = ".md"
$status    = "Handed back: in sync with en-US"
$xliffZh   = "afb6b190-c5eb-48b1-9fc3-5a122a05fee6.6b4462c086e7658335e0b54d2605e7de48340d22.zh-cn.xlf"
$xliffDe   = "afb6b190-c5eb-48b1-9fc3-5a122a05fee6.6b4462c086e7658335e0b54d2605e7de48340d22.de-de.xlf"

$hoDateZh  = "2016-08-24 12:47:18"
$hbDateZh  = "2016-08-24 12:47:35"
$hoDateDe  = "2016-08-24 12:47:23"
$hbDateDe  = "2016-08-24 12:47:42"
$genDate   = "2016-08-24 12:47:23"

$srcUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f0e2a9a8e0a7a7a4a02e6b3a3c9b8a9f1e4c2d7/e2e/afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md"
$zhUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6b4462c086e7658335e0b54d2605e7de48340d22/e2e/afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md"
$deUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6b4462c086e7658335e0b54d2605e7de48340d22/e2e/afb6b190-c5eb-48b1-9fc3-5a122a05fee6.md"

# ---------------------------------------------------------------------------
# Overview sheet (row 4): File Name | Path And Name | Extension | Publish URL
#                          | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overview.Range("A4").Value = $srcFile
$overview.Range("C4").Value = $ext
$overview.Range("E4").Value = $status
$overview.Range("F4").Value = $status
$overview.Range("G4").Value = $genDate
$overview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$overview.Hyperlinks.Add($overview.Range("B4"), $srcUrl, "", "", $srcPath)

# ---------------------------------------------------------------------------
# zh-cn sheet (row 4): Source File Name | File Extension | Status | Source Path
#  | Priority | Content Duplicate | Correspond Handoff File | Correspond Handoff Datetime
#  | Target File | Correspond Handback File | Correspond Handback DateTime
#  | Reference Tokens | To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$zhcn.Range("B4").Value = $ext
$zhcn.Range("C4").Value = $status
$zhcn.Range("D4").Value = "e2e"
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("F4").Value = "True"
$zhcn.Range("G4").Value = $xliffZh
$zhcn.Range("H4").Value = $hoDateZh
$zhcn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("J4").Value = $xliffZh
$zhcn.Range("K4").Value = $hbDateZh
$zhcn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("L4").Value = ""
$zhcn.Range("M4").Value = "True"
$zhcn.Range("N4").Value = ""
$zhcn.Range("O4").Value = "False"
$zhcn.Range("P4").Value = ""

$zhcn.Hyperlinks.Add($zhcn.Range("A4"), $srcUrl, "", "", $srcFile)
$zhcn.Hyperlinks.Add($zhcn.Range("I4"), $zhUrl, "", "", $srcFile)

# ---------------------------------------------------------------------------
# de-de sheet (row 4): same columns as zh-cn sheet
# ---------------------------------------------------------------------------
$dede.Range("B4").Value = $ext
$dede.Range("C4").Value = $status
$dede.Range("D4").Value = "e2e"
$dede.Range("E4").Value = "ht"
$dede.Range("F4").Value = "True"
$dede.Range("G4").Value = $xliffDe
$dede.Range("H4").Value = $hoDateDe
$dede.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("J4").Value = $xliffDe
$dede.Range("K4").Value = $hbDateDe
$dede.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("L4").Value = ""
$dede.Range("M4").Value = "True"
$dede.Range("N4").Value = ""
$dede.Range("O4").Value = "False"
$dede.Range("P4").Value = ""

$dede.Hyperlinks.Add($dede.Range("A4"), $srcUrl, "", "", $srcFile)
$dede.Hyperlinks.Add($dede.Range("I4"), $deUrl, "", "", $srcFile)

# ---------------------------------------------------------------------------
# Extend the tables to cover the newly added row.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Overview").ListObjects.Item(1).Resize($overview.Range("A1:G4"))
$wb.Worksheets.Item("zh-cn").ListObjects.Item(1).Resize($zhcn.Range("A1:P4"))
$wb.Worksheets.Item("de-de").ListObjects.Item(1).Resize($dede.Range("A1:P4"))
